$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Before state (tail of the document), 1-indexed paragraphs:
#   11: "2023.3．10"                         (no pPr)
#   12: "计划有变。" + "不知变不变。"           (pPr: rFonts hint=eastAsia)
#   13: ""                                    (pPr: rFonts hint=eastAsia)
#   14: "2023.3.23"                           (pPr: rFonts hint=eastAsia)
#   15: "没变。"                               (pPr: rFonts hint=eastAsia)
#   16: "今天学分支。" + bookmark "_GoBack"    (no pPr)
#
# After state (tail of the document):
#   11: "2023.3．10"                          (unchanged)
#   12: "计划有变。" + "不知变不变。"           (no pPr)            <- NEW paragraph
#   13: ""                                    (no pPr)            <- NEW paragraph
#   14: "2023.3.23"                           (no pPr)            <- NEW paragraph
#   15: "没变。"                               (no pPr)            <- NEW paragraph
#   16: "今天学分支。"                         (pPr: rFonts eastAsia) <- reuse of old 12
#   17: ""                                    (pPr: rFonts eastAsia) <- reuse of old 13 (unchanged)
#   18: "2023.3.23 " + "下午"                  (pPr: rFonts eastAsia) <- reuse of old 14
#   19: "笑死了，...技术"。 + bookmark          (no pPr)            <- reuse of old 16 (text changed)
#
# (old paragraph 15 "没变。" is no longer needed post-edit and is deleted)
#
# Paragraphs inserted via Range.InsertParagraphAfter()/Before() inherit the
# paragraph-mark formatting (pPr) of the paragraph whose Range owns the call,
# so inserting the four new "plain" paragraphs right after paragraph 11 (no
# pPr) keeps them free of the eastAsia paragraph-mark override, while the
# reused paragraphs 12/13/14 (which already carry that pPr) keep it as-is.
# ---------------------------------------------------------------------------

$anchor = $d.Paragraphs.Item(11)

# 1) Insert the four new "plain" paragraphs right after paragraph 11.
$anchor.Range.InsertParagraphAfter()
$p12 = $d.Paragraphs.Item(12)
$p12.Range.Text = "计划有变。不知变不变。"

$p12.Range.InsertParagraphAfter()
$p13 = $d.Paragraphs.Item(13)
# paragraph 13 stays empty

$p13.Range.InsertParagraphAfter()
$p14 = $d.Paragraphs.Item(14)
$p14.Range.Text = "2023.3.23"

$p14.Range.InsertParagraphAfter()
$p15 = $d.Paragraphs.Item(15)
$p15.Range.Text = "没变。"

# ---------------------------------------------------------------------------
# 2) The old paragraphs have all shifted down by four. What used to be
#    paragraph 12 ("计划有变。不知变不变。", pPr eastAsia) is now 16, the old
#    13 (empty, pPr eastAsia) is now 17, old 14 ("2023.3.23", pPr eastAsia)
#    is now 18, old 15 ("没变。", pPr eastAsia) is now 19, and old 16
#    ("今天学分支。" + bookmark, no pPr) is now 20.
# ---------------------------------------------------------------------------

# Reuse old-12 (now 16) -> "今天学分支。"
$p16 = $d.Paragraphs.Item(16)
$p16.Range.Text = "今天学分支。"

# old-13 (now 17) stays empty - nothing to do.

# Reuse old-14 (now 18) -> "2023.3.23 " + "下午" (two runs)
$p18 = $d.Paragraphs.Item(18)
$p18.Range.Text = "2023.3.23 "
$markPos = $p18.Range.End - 1
$d.Range($markPos, $markPos).InsertBefore("下午")

# old-15 (now 19, "没变。") is no longer needed - delete the whole paragraph.
$p19 = $d.Paragraphs.Item(19)
$p19.Range.Delete()

# Reuse old-16 (now 19 again, "今天学分支。" + bookmark) -> new closing text.
$p19b = $d.Paragraphs.Item(19)
$bookmarkName = "_GoBack"
$hasBookmark = $d.Bookmarks.Exists($bookmarkName)
$textEnd = $p19b.Range.End - 1
$textRange = $d.Range($p19b.Range.Start, $textEnd)
$textRange.Text = "笑死了，宿舍奇葩学姐被我们班男生集体吃瓜，老师：" + [char]34 + "你们在探讨什么技术" + [char]34 + "。"

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
